# Applies scheduled market-price / profit data refresh across all Leve sheets.
# Values below were recomputed upstream (e.g. from a market-data API) and
# are written verbatim into the currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$updates_ALC = @(
    @{ Cell = "H6"; Value = 37.5 }
    @{ Cell = "I6"; Value = 40 }
    @{ Cell = "J6"; Value = 35 }
    @{ Cell = "K6"; Value = 120 }
    @{ Cell = "L6"; Value = 105 }
    @{ Cell = "M6"; Value = -8 }
    @{ Cell = "N6"; Value = -329 }
    @{ Cell = "H74"; Value = 5670.2354 }
    @{ Cell = "I74"; Value = 4940.8 }
    @{ Cell = "K74"; Value = 4940.8 }
    @{ Cell = "M74"; Value = -4004.8 }
    @{ Cell = "H77"; Value = 5670.2354 }
    @{ Cell = "I77"; Value = 4940.8 }
    @{ Cell = "K77"; Value = 24704 }
    @{ Cell = "M77"; Value = -20024 }
    @{ Cell = "H88"; Value = 3668.9 }
    @{ Cell = "I88"; Value = 2446.5 }
    @{ Cell = "J88"; Value = 3974.5 }
    @{ Cell = "K88"; Value = 2446.5 }
    @{ Cell = "L88"; Value = 3974.5 }
    @{ Cell = "M88"; Value = -2040.5 }
    @{ Cell = "N88"; Value = -4786.5 }
    @{ Cell = "H91"; Value = 3668.9 }
    @{ Cell = "I91"; Value = 2446.5 }
    @{ Cell = "J91"; Value = 3974.5 }
    @{ Cell = "K91"; Value = 2446.5 }
    @{ Cell = "L91"; Value = 3974.5 }
    @{ Cell = "M91"; Value = -1042.5 }
    @{ Cell = "N91"; Value = -6782.5 }
    @{ Cell = "H131"; Value = 3247.5 }
    @{ Cell = "I131"; Value = 3247.5 }
    @{ Cell = "K131"; Value = 9742.5 }
    @{ Cell = "M131"; Value = -4702.5 }
    @{ Cell = "H132"; Value = 2774.1086 }
    @{ Cell = "I132"; Value = 2981.6052 }
    @{ Cell = "K132"; Value = 8944.8156 }
    @{ Cell = "M132"; Value = -6414.8156 }
    @{ Cell = "H137"; Value = 17385.148 }
    @{ Cell = "I137"; Value = 19884.842 }
    @{ Cell = "J137"; Value = 11448.375 }
    @{ Cell = "K137"; Value = 59654.526 }
    @{ Cell = "L137"; Value = 34345.125 }
    @{ Cell = "M137"; Value = -57104.526 }
    @{ Cell = "N137"; Value = -39445.125 }
)
foreach ($u in $updates_ALC) {
    $ws.Range($u.Cell).Value = $u.Value
}

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$updates_ARM = @(
    @{ Cell = "H4"; Value = 675.25 }
    @{ Cell = "J4"; Value = 400 }
    @{ Cell = "L4"; Value = 400 }
    @{ Cell = "H32"; Value = 3737.8413 }
    @{ Cell = "I32"; Value = 3704.6558 }
    @{ Cell = "K32"; Value = 3704.6558 }
    @{ Cell = "M32"; Value = -3417.6558 }
    @{ Cell = "H45"; Value = 4730.3335 }
    @{ Cell = "I45"; Value = 3500 }
    @{ Cell = "J45"; Value = 5345.5 }
    @{ Cell = "K45"; Value = 3500 }
    @{ Cell = "L45"; Value = 5345.5 }
    @{ Cell = "M45"; Value = -3123 }
    @{ Cell = "N45"; Value = -6099.5 }
    @{ Cell = "H46"; Value = 6418.6 }
    @{ Cell = "J46"; Value = 6418.6 }
    @{ Cell = "L46"; Value = 6418.6 }
    @{ Cell = "N46"; Value = -7056.6 }
    @{ Cell = "H97"; Value = 1976.6 }
    @{ Cell = "I97"; Value = 1943.3334 }
    @{ Cell = "J97"; Value = 2109.6667 }
    @{ Cell = "K97"; Value = 1943.3334 }
    @{ Cell = "L97"; Value = 2109.6667 }
    @{ Cell = "M97"; Value = -1447.3334 }
    @{ Cell = "N97"; Value = -3101.6667 }
    @{ Cell = "H126"; Value = 6565.3335 }
    @{ Cell = "I126"; Value = 6565.3335 }
    @{ Cell = "K126"; Value = 19696.0005 }
    @{ Cell = "M126"; Value = -17226.0005 }
    @{ Cell = "H132"; Value = 49272.234 }
    @{ Cell = "I132"; Value = 2613.1072 }
    @{ Cell = "K132"; Value = 7839.321599999999 }
    @{ Cell = "M132"; Value = -5309.321599999999 }
    @{ Cell = "H135"; Value = 113168 }
    @{ Cell = "J135"; Value = 113168 }
    @{ Cell = "L135"; Value = 113168 }
    @{ Cell = "N135"; Value = -123308 }
)
foreach ($u in $updates_ARM) {
    $ws.Range($u.Cell).Value = $u.Value
}

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$updates_BSM = @(
    @{ Cell = "H20"; Value = 1681.2 }
    @{ Cell = "I20"; Value = 1727.5 }
    @{ Cell = "J20"; Value = 1628.2858 }
    @{ Cell = "K20"; Value = 1727.5 }
    @{ Cell = "L20"; Value = 1628.2858 }
    @{ Cell = "M20"; Value = -1480.5 }
    @{ Cell = "N20"; Value = -2122.2858 }
    @{ Cell = "H40"; Value = 493493 }
    @{ Cell = "J40"; Value = 493493 }
    @{ Cell = "L40"; Value = 493493 }
    @{ Cell = "N40"; Value = -494023 }
    @{ Cell = "H86"; Value = 15311.546 }
    @{ Cell = "I86"; Value = 9939.166999999999 }
    @{ Cell = "J86"; Value = 21758.4 }
    @{ Cell = "K86"; Value = 9939.166999999999 }
    @{ Cell = "L86"; Value = 21758.4 }
    @{ Cell = "M86"; Value = -8816.166999999999 }
    @{ Cell = "N86"; Value = -24004.4 }
    @{ Cell = "H89"; Value = 15311.546 }
    @{ Cell = "I89"; Value = 9939.166999999999 }
    @{ Cell = "J89"; Value = 21758.4 }
    @{ Cell = "K89"; Value = 49695.835 }
    @{ Cell = "L89"; Value = 108792 }
    @{ Cell = "M89"; Value = -44079.835 }
    @{ Cell = "N89"; Value = -120024 }
    @{ Cell = "H94"; Value = 1267.0741 }
    @{ Cell = "J94"; Value = 2142.7273 }
    @{ Cell = "L94"; Value = 2142.7273 }
    @{ Cell = "N94"; Value = -3044.7273 }
    @{ Cell = "H107"; Value = 2129.182 }
    @{ Cell = "I107"; Value = 1908.9667 }
    @{ Cell = "K107"; Value = 1908.9667 }
    @{ Cell = "M107"; Value = 11.03330000000005 }
    @{ Cell = "H134"; Value = 1106 }
    @{ Cell = "I134"; Value = 1106 }
    @{ Cell = "K134"; Value = 3318 }
    @{ Cell = "M134"; Value = -783 }
)
foreach ($u in $updates_BSM) {
    $ws.Range($u.Cell).Value = $u.Value
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$updates_CRP = @(
    @{ Cell = "H99"; Value = 5018.6665 }
    @{ Cell = "I99"; Value = 4324.0713 }
    @{ Cell = "J99"; Value = 7449.75 }
    @{ Cell = "K99"; Value = 4324.0713 }
    @{ Cell = "L99"; Value = 7449.75 }
    @{ Cell = "M99"; Value = -2826.0713 }
    @{ Cell = "N99"; Value = -10445.75 }
    @{ Cell = "H126"; Value = 5018.6665 }
    @{ Cell = "I126"; Value = 4324.0713 }
    @{ Cell = "J126"; Value = 7449.75 }
    @{ Cell = "K126"; Value = 12972.2139 }
    @{ Cell = "L126"; Value = 22349.25 }
    @{ Cell = "M126"; Value = -10502.2139 }
    @{ Cell = "N126"; Value = -27289.25 }
    @{ Cell = "H132"; Value = 4321.875 }
    @{ Cell = "I132"; Value = 4105.5186 }
    @{ Cell = "J132"; Value = 5490.2 }
    @{ Cell = "K132"; Value = 12316.5558 }
    @{ Cell = "L132"; Value = 16470.6 }
    @{ Cell = "M132"; Value = -9786.555800000002 }
    @{ Cell = "N132"; Value = -21530.6 }
)
foreach ($u in $updates_CRP) {
    $ws.Range($u.Cell).Value = $u.Value
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$updates_CUL = @(
    @{ Cell = "H17"; Value = 157 }
    @{ Cell = "I17"; Value = 157 }
    @{ Cell = "K17"; Value = 471 }
    @{ Cell = "M17"; Value = -302 }
    @{ Cell = "H131"; Value = 1167.6666 }
    @{ Cell = "J131"; Value = 1627 }
    @{ Cell = "L131"; Value = 4881 }
    @{ Cell = "N131"; Value = -14961 }
    @{ Cell = "H137"; Value = 3815.8333 }
    @{ Cell = "I137"; Value = 1030 }
    @{ Cell = "J137"; Value = 4373 }
    @{ Cell = "K137"; Value = 3090 }
    @{ Cell = "L137"; Value = 13119 }
    @{ Cell = "M137"; Value = 2010 }
    @{ Cell = "N137"; Value = -23319 }
)
foreach ($u in $updates_CUL) {
    $ws.Range($u.Cell).Value = $u.Value
}

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$updates_GSM = @(
    @{ Cell = "H70"; Value = 6278.8 }
    @{ Cell = "I70"; Value = 5625.857 }
    @{ Cell = "K70"; Value = 5625.857 }
    @{ Cell = "M70"; Value = -5355.857 }
    @{ Cell = "H73"; Value = 6278.8 }
    @{ Cell = "I73"; Value = 5625.857 }
    @{ Cell = "K73"; Value = 5625.857 }
    @{ Cell = "M73"; Value = -4689.857 }
    @{ Cell = "H97"; Value = 832.5 }
    @{ Cell = "I97"; Value = 850.55554 }
    @{ Cell = "J97"; Value = 800 }
    @{ Cell = "K97"; Value = 850.55554 }
    @{ Cell = "L97"; Value = 800 }
    @{ Cell = "M97"; Value = -354.55554 }
    @{ Cell = "N97"; Value = -1792 }
    @{ Cell = "H102"; Value = 4124.2905 }
    @{ Cell = "I102"; Value = 2146.1738 }
    @{ Cell = "K102"; Value = 2146.1738 }
    @{ Cell = "M102"; Value = -524.1738 }
    @{ Cell = "H132"; Value = 1757 }
    @{ Cell = "I132"; Value = 0 }
    @{ Cell = "J132"; Value = 1757 }
    @{ Cell = "K132"; Value = 0 }
    @{ Cell = "L132"; Value = 5271 }
    @{ Cell = "N132"; Value = -10331 }
)
foreach ($u in $updates_GSM) {
    $ws.Range($u.Cell).Value = $u.Value
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$updates_LTW = @(
    @{ Cell = "H16"; Value = 1054.8148 }
    @{ Cell = "I16"; Value = 979.2 }
    @{ Cell = "K16"; Value = 979.2 }
    @{ Cell = "M16"; Value = -809.2 }
    @{ Cell = "H20"; Value = 846 }
    @{ Cell = "I20"; Value = 652 }
    @{ Cell = "K20"; Value = 652 }
    @{ Cell = "M20"; Value = -426 }
    @{ Cell = "H22"; Value = 1325.8695 }
    @{ Cell = "I22"; Value = 1291.5714 }
    @{ Cell = "J22"; Value = 1379.2222 }
    @{ Cell = "K22"; Value = 1291.5714 }
    @{ Cell = "L22"; Value = 1379.2222 }
    @{ Cell = "M22"; Value = -996.5714 }
    @{ Cell = "N22"; Value = -1969.2222 }
    @{ Cell = "H27"; Value = 1325.8695 }
    @{ Cell = "I27"; Value = 1291.5714 }
    @{ Cell = "J27"; Value = 1379.2222 }
    @{ Cell = "K27"; Value = 1291.5714 }
    @{ Cell = "L27"; Value = 1379.2222 }
    @{ Cell = "M27"; Value = -1184.5714 }
    @{ Cell = "N27"; Value = -1593.2222 }
    @{ Cell = "H55"; Value = 756.8889 }
    @{ Cell = "I55"; Value = 477.35294 }
    @{ Cell = "J55"; Value = 1232.1 }
    @{ Cell = "K55"; Value = 477.35294 }
    @{ Cell = "L55"; Value = 1232.1 }
    @{ Cell = "M55"; Value = -304.35294 }
    @{ Cell = "N55"; Value = -1578.1 }
    @{ Cell = "H93"; Value = 4256.125 }
    @{ Cell = "I93"; Value = 3463.75 }
    @{ Cell = "J93"; Value = 5048.5 }
    @{ Cell = "K93"; Value = 3463.75 }
    @{ Cell = "L93"; Value = 5048.5 }
    @{ Cell = "M93"; Value = -2215.75 }
    @{ Cell = "N93"; Value = -7544.5 }
    @{ Cell = "H132"; Value = 2395.6875 }
    @{ Cell = "I132"; Value = 2252.8076 }
    @{ Cell = "J132"; Value = 3014.8333 }
    @{ Cell = "K132"; Value = 6758.4228 }
    @{ Cell = "L132"; Value = 9044.499899999999 }
    @{ Cell = "M132"; Value = -4228.4228 }
    @{ Cell = "N132"; Value = -14104.4999 }
    @{ Cell = "H133"; Value = 68122.5 }
    @{ Cell = "J133"; Value = 68122.5 }
    @{ Cell = "L133"; Value = 68122.5 }
    @{ Cell = "N133"; Value = -73182.5 }
    @{ Cell = "H136"; Value = 2260.6667 }
    @{ Cell = "I136"; Value = 2237.9678 }
    @{ Cell = "K136"; Value = 6713.903399999999 }
    @{ Cell = "M136"; Value = -4163.903399999999 }
)
foreach ($u in $updates_LTW) {
    $ws.Range($u.Cell).Value = $u.Value
}

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$updates_WVR = @(
    @{ Cell = "H132"; Value = 1033.6666 }
    @{ Cell = "I132"; Value = 1033.6666 }
    @{ Cell = "K132"; Value = 3100.9998 }
    @{ Cell = "M132"; Value = -570.9998000000001 }
    @{ Cell = "H133"; Value = 85022.60000000001 }
    @{ Cell = "J133"; Value = 85022.60000000001 }
    @{ Cell = "L133"; Value = 85022.60000000001 }
    @{ Cell = "N133"; Value = -95142.60000000001 }
    @{ Cell = "H136"; Value = 3054.365 }
    @{ Cell = "I136"; Value = 1999.4565 }
    @{ Cell = "K136"; Value = 5998.3695 }
    @{ Cell = "M136"; Value = -3448.3695 }
)
foreach ($u in $updates_WVR) {
    $ws.Range($u.Cell).Value = $u.Value
}

# ---- Special case: ARM row 4 previously had no LeveProfitHQ (N4) value; now populated ----
$wsARM = $wb.Worksheets.Item("ARM")
$wsARM.Range("N4").Value = -632

# ---- Special case: GSM row 132 LeveProfitNQ (M132) column removed (HQ-only leve now) ----
$wsGSM = $wb.Worksheets.Item("GSM")
$wsGSM.Range("M132").ClearContents()
